$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.832.78"
$ws.Range("E2").Value = "  +5.93%  "

$ws.Range("D3").Value = "3.475.96"
$ws.Range("E3").Value = "  +3.85%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'412.48"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").Value = "'128.91"
$ws.Range("E6").Value = "  +12.49%  "

$ws.Range("D7").Value = "3.467.51"
$ws.Range("E7").Value = "  +3.85%  "

$ws.Range("D8").Value = "'0.601"
$ws.Range("E8").Value = "  +2.04%  "

$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").Value = "'0.703"
$ws.Range("E10").Value = "  +9.50%  "

$ws.Range("E11").Value = "  +29.74%  "

$ws.Range("D12").Value = "'43.38"
$ws.Range("E12").Value = "  +7.04%  "

$ws.Range("E13").Value = "  -0.54%  "

$ws.Range("D14").Value = "4.025.69"
$ws.Range("E14").Value = "  +3.86%  "

$ws.Range("E15").Value = "  +2.62%  "

$ws.Range("D16").Value = "'20.26"
$ws.Range("E16").Value = "  +3.59%  "

$ws.Range("D17").Value = "3.450.22"
$ws.Range("E17").Value = "  +3.12%  "

$ws.Range("D18").Value = "62.669.83"
$ws.Range("E18").Value = "  +6.06%  "

$ws.Range("D19").Value = "'1.04"
$ws.Range("E19").Value = "  -0.47%  "

$ws.Range("D20").Value = "'11.11"
$ws.Range("E20").Value = "  +2.32%  "

$ws.Range("E21").Value = "  +23.10%  "

$ws.Range("D22").Value = "'3.35"
$ws.Range("E22").Value = "  -0.69%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'82.41"
$ws.Range("E23").Value = "  +8.56%  "

$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "'13.21"
$ws.Range("E24").Value = "  -0.30%  "

$ws.Range("D25").Value = "'314.48"
$ws.Range("E25").Value = "  +3.16%  "

$ws.Range("D26").Value = "'3.18"
$ws.Range("E26").Value = "  -0.59%  "

$ws.Range("E27").Value = "  +5.61%  "

$ws.Range("D28").Value = "'8.20"
$ws.Range("E28").Value = "  +3.22%  "

$ws.Range("D29").Value = "'7.80"
$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("D30").Value = "'0.182"
$ws.Range("E30").Value = "  +3.27%  "

$ws.Range("B31").Value = "LEO"
$ws.Range("C31").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D31").Value = "'4.37"
$ws.Range("E31").Value = "  -2.80%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.120"
$ws.Range("E32").Value = "  +1.55%  "

$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'12.15"
$ws.Range("E33").Value = "  +4.89%  "

$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "'44.39"
$ws.Range("E34").Value = "  +9.64%  "

$ws.Range("E35").Value = "  +25.79%  "

$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("E37").Value = "  -5.63%  "

$ws.Range("D38").Value = "'52.66"
$ws.Range("E38").Value = "  +1.28%  "

$ws.Range("D39").Value = "'3.59"
$ws.Range("E39").Value = "  +5.61%  "

$ws.Range("E40").Value = "  -0.26%  "

$ws.Range("D41").Value = "'3.03"
$ws.Range("E41").Value = "  -3.64%  "

$ws.Range("E42").Value = "  +3.10%  "

$ws.Range("E43").Value = "  +2.53%  "

$ws.Range("D44").Value = "'137.63"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").Value = "'17.81"
$ws.Range("E45").Value = "  +2.64%  "

$ws.Range("D46").Value = "'0.288"
$ws.Range("E46").Value = "  +2.16%  "

$ws.Range("D47").Value = "'3.98"
$ws.Range("E47").Value = "  -0.81%  "

$ws.Range("D48").Value = "'2.27"
$ws.Range("E48").Value = "  +0.97%  "

$ws.Range("E49").Value = "  -0.35%  "

$ws.Range("D50").Value = "2.222.96"
$ws.Range("E50").Value = "  +0.44%  "

$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").Value = "'2.38"
$ws.Range("E51").Value = "  -1.10%  "
